$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1099-C")

# Rename the "Is Corrected" header to "Is Corrected Form of 1099"
$ws.Range("X1").Value = "Is Corrected Form of 1099"

# Rows 9-13 (X9:X13) flip from "Yes" to "No"
$ws.Range("X9:X13").Value = "No"

# Match the author's final selection on the sheet
$ws.Range("X9:X13").Select() | Out-Null
